$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.053.42"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.829.61"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "240.73"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "0.6217"
$ws.Range("E6").Value = "  -6.40%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "44.45"
$ws.Range("E8").Value = "  +5.81%  "
$ws.Range("D9").Value = "0.07378"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "0.2923"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").Value = "22.67"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "1.831.06"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "4.958"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").Value = "0.6628"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "82.07"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "0.000009097"
$ws.Range("E17").Value = "  +8.80%  "
$ws.Range("D18").Value = "6.013"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "29.051.36"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "2.079.19"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "225.47"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "7.173"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "159.50"
$ws.Range("D27").Value = "8.422"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").Value = "0.1357"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("D29").Value = "17.81"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "1.494"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "4.054"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "4.028"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "1.202"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "0.05243"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").Value = "1.838"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "1.152"
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("D38").Value = "2.646"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").Value = "1.283.94"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "2.749"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "0.01780"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "6.326"
$ws.Range("E42").Value = "  +6.67%  "
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "101.73"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "1.976.82"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").Value = "63.73"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "0.00000000119"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("D51").Value = "0.3966"
$ws.Range("E51").Value = "  -1.47%  "
